$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 37 (first new row of tracked work)
$ws.Range("A37").NumberFormat = "d-mmm"
$ws.Range("A37").Value = (Get-Date -Year 2020 -Month 3 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B37").Value = 11
$ws.Range("C37").Value = "abandon de la library dropzone page ajout theme/emission podcast fonctionelle"

# Add row 38 (second new row of tracked work)
$ws.Range("A38").NumberFormat = "d-mmm"
$ws.Range("A38").Value = (Get-Date -Year 2020 -Month 3 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B38").Value = 8
$ws.Range("C38").Value = "page modifier en cours de construction modifier theme et emission fonctionelle"

# Update the two existing "DropFile" -> "Dropzone" descriptions
$ws.Range("C23").Value = "etude Dropzone"
$ws.Range("C24").Value = "outil Dropzone maitrisé et modilation de model pour future implémentation"

# Add row 39 (third new row of tracked work)
$ws.Range("A39").NumberFormat = "d-mmm"
$ws.Range("A39").Value = (Get-Date -Year 2020 -Month 3 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B39").Value = 3
$ws.Range("C39").Value = "page modifier fonctionelle site fonctionelle v1"

# Recalculate so F6 (SUM(B7:B39)) reflects the new total
$excel.Calculate()

# Move/update the active selection to C39 as in the edited workbook
$ws.Range("C39").Select()
